# Generate Report for Handoff
#
# The localization-status report is regenerated: the two tracked files swap
# row order (the "29f496bf" file now reports first, the "014c6cdf" file now
# reports second) and the "014c6cdf" file's status flips from "In Translation"
# to "Ready for handoff" with refreshed handoff timestamps. This touches the
# "Overview" sheet plus the per-locale "zh-cn" and "de-de" sheets. Existing
# hyperlink relationships (r:id -> external URL) are left alone; only the
# cell text (and therefore each hyperlink's auto-synced display text) changes.

$wb = $excel.ActiveWorkbook

$mdA = "014c6cdf-958b-4bc2-b78a-053e8c86b8b8.md"
$mdB = "29f496bf-9dac-48a4-b9dc-5f582194e5ec.md"

$statusOld = "In Translation"
$statusNew = "Ready for handoff"

# ---- Overview sheet ------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = $mdB
$ov.Range("B2").Value = $statusOld
$ov.Range("C2").Value = $statusOld
$ov.Range("D2").Value = "2016-03-23 04:20:56"

$ov.Range("A3").Value = $mdA
$ov.Range("B3").Value = $statusNew
$ov.Range("C3").Value = $statusNew
$ov.Range("D3").Value = "2016-03-23 04:22:07"

# ---- zh-cn sheet ----------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = $mdB
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = $statusOld
$zh.Range("D2").Value = "29f496bf-9dac-48a4-b9dc-5f582194e5ec.19ede27134ecbd1a9e7894f4129a6a54e30b2928.zh-cn.xlf"
$zh.Range("E2").Value = "2016-03-23 04:20:53"
$zh.Range("H2").Value = "0001-01-01 00:00:00"
$zh.Range("J2").Value = "Include"

$zh.Range("A3").Value = $mdA
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = $statusNew
$zh.Range("D3").Value = "014c6cdf-958b-4bc2-b78a-053e8c86b8b8.e44e71d4f0489edd6755148b97b69e11f7257c4a.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-23 04:22:02"
$zh.Range("H3").Value = "0001-01-01 00:00:00"
$zh.Range("J3").Value = "Include"

# ---- de-de sheet ----------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = $mdB
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = $statusOld
$de.Range("D2").Value = "29f496bf-9dac-48a4-b9dc-5f582194e5ec.19ede27134ecbd1a9e7894f4129a6a54e30b2928.de-de.xlf"
$de.Range("E2").Value = "2016-03-23 04:20:56"
$de.Range("H2").Value = "0001-01-01 00:00:00"
$de.Range("J2").Value = "Include"

$de.Range("A3").Value = $mdA
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = $statusNew
$de.Range("D3").Value = "014c6cdf-958b-4bc2-b78a-053e8c86b8b8.e44e71d4f0489edd6755148b97b69e11f7257c4a.de-de.xlf"
$de.Range("E3").Value = "2016-03-23 04:22:07"
$de.Range("H3").Value = "0001-01-01 00:00:00"
$de.Range("J3").Value = "Include"
